# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to match the refreshed cryptos data from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column cells we touch so that purely
# numeric-looking strings (e.g. "556.60") are kept as text, matching the
# original inlineStr/text cells instead of being parsed into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.774.15"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.052.48"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.60"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.90"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.050.94"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  -11.06%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.12"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.550.08"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.804.62"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.045.69"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.62"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.42"
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.51"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.51"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.10"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.27"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.19"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.67"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0408"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "441.61"
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0814"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.023.39"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.31"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.271"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.63"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.71"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0513"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  +0.31%  "
